$d = $word.ActiveDocument

function Replace-ParagraphXml($findText, $innerXml) {
    $r = $d.Content
    $found = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find paragraph text: $findText"
    }
    $para = $r.Paragraphs(1)
    $target = $para.Range
    $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($xml) | Out-Null
}

# --- Part 1 --------------------------------------------------------------
# "LUẬN VĂN TỐT NGHIỆP ĐẠI HỌC" is split into two runs with an empty
# "_GoBack" bookmark inserted between "NGHI" and "ỆP".
$fullTitle = "LUẬN VĂN TỐT NGHIỆP ĐẠI HỌC"
$titlePrefix = "LUẬN VĂN TỐT NGHI"

$titleRange = $d.Content
$found = $titleRange.Find.Execute($fullTitle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find title text"
}
$splitPos = $titleRange.Start + $titlePrefix.Length
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- Part 2 ----------------------------------------------------------------
# The three paragraphs that make up the thesis subtitle are reformatted
# (centered, size 32, tabs removed) and their text content changed.
$pPrCommon = '<w:pPr><w:spacing w:after="0"/><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr>'
$rPrCommon = '<w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>'

# XÂY DỰNG HỆ THỐNG THÔNG TIN ĐỊA LÝ
$xmlA = '<w:p>' + $pPrCommon + '<w:r>' + $rPrCommon + '<w:t>XÂY DỰNG HỆ THỐNG THÔNG TIN ĐỊA LÝ</w:t></w:r></w:p>'
Replace-ParagraphXml "XÂY DỰNG HỆ THỐNG THÔNG TIN ĐỊA LÝ" $xmlA

# VỀ HẠ TẦNG GIAO THÔNG BỘ THÀNH PHỐ CẦN THƠ (two runs)
$xmlB = '<w:p>' + $pPrCommon + '<w:r>' + $rPrCommon + '<w:t>VỀ HẠ TẦNG GIAO THÔNG</w:t></w:r><w:r>' + $rPrCommon + '<w:t xml:space="preserve"> BỘ THÀNH PHỐ CẦN THƠ</w:t></w:r></w:p>'
Replace-ParagraphXml "VỀ HẠ TẦNG GIAO THÔNG" $xmlB

# "Phân hệ giao thông bộ" paragraph becomes empty (text removed)
$xmlC = '<w:p>' + $pPrCommon + '</w:p>'
Replace-ParagraphXml "Phân hệ giao thông bộ" $xmlC
